$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.043711185455322
$ws.Range("B1").Value = 2.084349393844604
$ws.Range("C1").Value = 2.593115568161011
$ws.Range("D1").Value = 2.265393733978271
$ws.Range("E1").Value = 2.004603624343872
